$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description text in B2: remove "RME/" from the "20% S/LFM+CDM/RME/H:1" line
$ws.Range("B2").Replace("20% S/LFM+CDM/RME/H:1", "20% S/LFM+CDM/H:1")

# Widen column B and wrap text / grow row 2 to fit the (now shorter) text block
$ws.Columns("B").ColumnWidth = 29.6640625
$ws.Range("B2").WrapText = $true
$ws.Rows(2).RowHeight = 96

# Move active selection to C2
$ws.Range("C2").Select()
